# Apply "Natmi following Dr Hou advice" recalculated values to the LR-pairs sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "FAPs"
$ws.Cells.Item(2, 2).Value = "Wnt5a"
$ws.Cells.Item(2, 3).Value = "Fzd5"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 4.619088000000001
$ws.Cells.Item(2, 8).Value = 13.857264
$ws.Cells.Item(2, 9).Value = 1
$ws.Cells.Item(2, 10).Value = 1
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.5
$ws.Cells.Item(2, 13).Value = 2.170377
$ws.Cells.Item(2, 14).Value = 4.340754
$ws.Cells.Item(2, 15).Value = 0.1015511790371702
$ws.Cells.Item(2, 16).Value = 0.07285982038608425
$ws.Cells.Item(2, 17).Value = 10.025162356176
$ws.Cells.Item(2, 18).Value = 60.150974137056
$ws.Cells.Item(2, 19).Value = 0.1015511790371702
$ws.Cells.Item(2, 20).Value = 0.07285982038608425

# Row 3
$ws.Cells.Item(3, 1).Value = "FAPs"
$ws.Cells.Item(3, 2).Value = "Wnt5a"
$ws.Cells.Item(3, 3).Value = "Fzd5"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 4.619088000000001
$ws.Cells.Item(3, 8).Value = 13.857264
$ws.Cells.Item(3, 9).Value = 1
$ws.Cells.Item(3, 10).Value = 1
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 5.061974333333333
$ws.Cells.Item(3, 14).Value = 15.185923
$ws.Cells.Item(3, 15).Value = 0.2368480046581279
$ws.Cells.Item(3, 16).Value = 0.2548966428820674
$ws.Cells.Item(3, 17).Value = 23.381704899408
$ws.Cells.Item(3, 18).Value = 210.435344094672
$ws.Cells.Item(3, 19).Value = 0.2368480046581279
$ws.Cells.Item(3, 20).Value = 0.2548966428820674

# Row 4
$ws.Cells.Item(4, 1).Value = "FAPs"
$ws.Cells.Item(4, 2).Value = "Wnt5a"
$ws.Cells.Item(4, 3).Value = "Fzd5"
$ws.Cells.Item(4, 4).Value = "M1"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 4.619088000000001
$ws.Cells.Item(4, 8).Value = 13.857264
$ws.Cells.Item(4, 9).Value = 1
$ws.Cells.Item(4, 10).Value = 1
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 3.815520666666667
$ws.Cells.Item(4, 14).Value = 11.446562
$ws.Cells.Item(4, 15).Value = 0.1785268745202745
$ws.Cells.Item(4, 16).Value = 0.1921312406457904
$ws.Cells.Item(4, 17).Value = 17.624225725152
$ws.Cells.Item(4, 18).Value = 158.618031526368
$ws.Cells.Item(4, 19).Value = 0.1785268745202745
$ws.Cells.Item(4, 20).Value = 0.1921312406457904

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Wnt5a"
$ws.Cells.Item(5, 3).Value = "Fzd5"
$ws.Cells.Item(5, 4).Value = "M2"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 4.619088000000001
$ws.Cells.Item(5, 8).Value = 13.857264
$ws.Cells.Item(5, 9).Value = 1
$ws.Cells.Item(5, 10).Value = 1
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 4.666218666666667
$ws.Cells.Item(5, 14).Value = 13.998656
$ws.Cells.Item(5, 15).Value = 0.2183307357409577
$ws.Cells.Item(5, 16).Value = 0.2349682939430755
$ws.Cells.Item(5, 17).Value = 21.553674648576
$ws.Cells.Item(5, 18).Value = 193.983071837184
$ws.Cells.Item(5, 19).Value = 0.2183307357409577
$ws.Cells.Item(5, 20).Value = 0.2349682939430755

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Wnt5a"
$ws.Cells.Item(6, 3).Value = "Fzd5"
$ws.Cells.Item(6, 4).Value = "Neutro"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 4.619088000000001
$ws.Cells.Item(6, 8).Value = 13.857264
$ws.Cells.Item(6, 9).Value = 1
$ws.Cells.Item(6, 10).Value = 1
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 3.288577
$ws.Cells.Item(6, 14).Value = 9.865731
$ws.Cells.Item(6, 15).Value = 0.1538713650690733
$ws.Cells.Item(6, 16).Value = 0.1655968959856798
$ws.Cells.Item(6, 17).Value = 15.190226557776
$ws.Cells.Item(6, 18).Value = 136.712039019984
$ws.Cells.Item(6, 19).Value = 0.1538713650690733
$ws.Cells.Item(6, 20).Value = 0.1655968959856798

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Wnt5a"
$ws.Cells.Item(7, 3).Value = "Fzd5"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 4.619088000000001
$ws.Cells.Item(7, 8).Value = 13.857264
$ws.Cells.Item(7, 9).Value = 1
$ws.Cells.Item(7, 10).Value = 1
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 2.3695805
$ws.Cells.Item(7, 14).Value = 4.739161
$ws.Cells.Item(7, 15).Value = 0.1108718409743963
$ws.Cells.Item(7, 16).Value = 0.07954710615730251
$ws.Cells.Item(7, 17).Value = 10.945300852584
$ws.Cells.Item(7, 18).Value = 65.671805115504
$ws.Cells.Item(7, 19).Value = 0.1108718409743963
$ws.Cells.Item(7, 20).Value = 0.07954710615730251
